# Update the sample laptop dataset (chore(data): update sample laptop dataset)
# - Replace the placeholder "rtx3050" gpu string with proper ram_gb/storage_gb
#   numeric columns plus real GPU model names.
# - Refresh several rows' RAM/storage/GPU/screen/weight/price figures.
# - Add RAM/storage/GPU/screen/weight/price data for the final row (Legion 5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = 7900
$ws.Cells.Item(2,5).Value = 8
$ws.Cells.Item(2,6).Value = 512
$ws.Cells.Item(2,8).Value = 15

# Row 3
$ws.Cells.Item(3,6).Value = 512
$ws.Cells.Item(3,7).Value = "Integrated Graphics"
$ws.Cells.Item(3,9).Value = 1

# Row 4
$ws.Cells.Item(4,6).Value = 512
$ws.Cells.Item(4,7).Value = "Integrated Graphics"
$ws.Cells.Item(4,8).Value = 15
$ws.Cells.Item(4,9).Value = 1

# Row 5
$ws.Cells.Item(5,5).Value = 8
$ws.Cells.Item(5,6).Value = 512
$ws.Cells.Item(5,7).Value = "Integrated Graphics"
$ws.Cells.Item(5,8).Value = 14
$ws.Cells.Item(5,9).Value = 1

# Row 6
$ws.Cells.Item(6,5).Value = 8
$ws.Cells.Item(6,6).Value = 512
$ws.Cells.Item(6,7).Value = "NVIDIA GeForce RTX 3050"
$ws.Cells.Item(6,8).Value = 15
$ws.Cells.Item(6,9).Value = 2

# Row 7
$ws.Cells.Item(7,5).Value = 8
$ws.Cells.Item(7,6).Value = 256
$ws.Cells.Item(7,7).Value = "Apple M1 GPU"
$ws.Cells.Item(7,8).Value = 13
$ws.Cells.Item(7,9).Value = 1
$ws.Cells.Item(7,10).Value = 14500000

# Row 8
$ws.Cells.Item(8,5).Value = 16
$ws.Cells.Item(8,6).Value = 512
$ws.Cells.Item(8,7).Value = "NVIDIA GeForce RTX 4050"
$ws.Cells.Item(8,8).Value = 15
$ws.Cells.Item(8,9).Value = 2
$ws.Cells.Item(8,10).Value = 15500000

# Row 9
$ws.Cells.Item(9,5).Value = 8
$ws.Cells.Item(9,6).Value = 512
$ws.Cells.Item(9,7).Value = "Integrated Graphics"
$ws.Cells.Item(9,8).Value = 14
$ws.Cells.Item(9,9).Value = 1
$ws.Cells.Item(9,10).Value = 12500000

# Row 10
$ws.Cells.Item(10,5).Value = 8
$ws.Cells.Item(10,6).Value = 512
$ws.Cells.Item(10,7).Value = "NVIDIA GeForce RTX 3050"
$ws.Cells.Item(10,8).Value = 14
$ws.Cells.Item(10,9).Value = 1
$ws.Cells.Item(10,10).Value = 11800000

# Row 11
$ws.Cells.Item(11,5).Value = 8
$ws.Cells.Item(11,6).Value = 512
$ws.Cells.Item(11,7).Value = "Integrated Graphics"
$ws.Cells.Item(11,8).Value = 14
$ws.Cells.Item(11,9).Value = 1
$ws.Cells.Item(11,10).Value = 13200000

# Row 12
$ws.Cells.Item(12,5).Value = 8
$ws.Cells.Item(12,6).Value = 512
$ws.Cells.Item(12,7).Value = "Intel Iris Xe Graphics"
$ws.Cells.Item(12,8).Value = 13
$ws.Cells.Item(12,9).Value = 1
$ws.Cells.Item(12,10).Value = 19000000

# Row 13
$ws.Cells.Item(13,5).Value = 16
$ws.Cells.Item(13,6).Value = 512
$ws.Cells.Item(13,7).Value = "NVIDIA GeForce RTX 4060"
$ws.Cells.Item(13,8).Value = 15
$ws.Cells.Item(13,9).Value = 2
$ws.Cells.Item(13,10).Value = 13000000

# Row 14
$ws.Cells.Item(14,5).Value = 16
$ws.Cells.Item(14,6).Value = 512
$ws.Cells.Item(14,7).Value = "NVIDIA GeForce RTX 3050"
$ws.Cells.Item(14,8).Value = 15
$ws.Cells.Item(14,9).Value = 2
$ws.Cells.Item(14,10).Value = 17000000

# Row 15
$ws.Cells.Item(15,5).Value = 8
$ws.Cells.Item(15,6).Value = 512
$ws.Cells.Item(15,7).Value = "Integrated Graphics"
$ws.Cells.Item(15,8).Value = 14
$ws.Cells.Item(15,9).Value = 1
$ws.Cells.Item(15,10).Value = 18500000

# Row 16
$ws.Cells.Item(16,5).Value = 8
$ws.Cells.Item(16,6).Value = 512
$ws.Cells.Item(16,7).Value = "NVIDIA GeForce RTX 4060"
$ws.Cells.Item(16,8).Value = 14
$ws.Cells.Item(16,9).Value = 1
$ws.Cells.Item(16,10).Value = 16500000

# Row 17
$ws.Cells.Item(17,5).Value = 16
$ws.Cells.Item(17,6).Value = 512
$ws.Cells.Item(17,7).Value = "NVIDIA GeForce RTX 4050"
$ws.Cells.Item(17,8).Value = 15
$ws.Cells.Item(17,9).Value = 2
$ws.Cells.Item(17,10).Value = 23000000

# Row 18
$ws.Cells.Item(18,5).Value = 16
$ws.Cells.Item(18,6).Value = 512
$ws.Cells.Item(18,7).Value = "Integrated Graphics"
$ws.Cells.Item(18,8).Value = 15
$ws.Cells.Item(18,9).Value = 2
$ws.Cells.Item(18,10).Value = 19500000

# Row 19
$ws.Cells.Item(19,5).Value = 8
$ws.Cells.Item(19,6).Value = 1024
$ws.Cells.Item(19,7).Value = "NVIDIA GeForce RTX 4070"
$ws.Cells.Item(19,8).Value = 13
$ws.Cells.Item(19,9).Value = 1
$ws.Cells.Item(19,10).Value = 20000000

# Row 20
$ws.Cells.Item(20,5).Value = 16
$ws.Cells.Item(20,6).Value = 1024
$ws.Cells.Item(20,7).Value = "NVIDIA GeForce RTX 4060"
$ws.Cells.Item(20,8).Value = 15
$ws.Cells.Item(20,9).Value = 2
$ws.Cells.Item(20,10).Value = 28000000

# Row 21
$ws.Cells.Item(21,5).Value = 16
$ws.Cells.Item(21,6).Value = 1024
$ws.Cells.Item(21,7).Value = "NVIDIA GeForce RTX 4060"
$ws.Cells.Item(21,8).Value = 15
$ws.Cells.Item(21,9).Value = 2
$ws.Cells.Item(21,10).Value = 24000000


# Widen the gpu column to fit the longer GPU names now stored in it.
$ws.Columns.Item(7).ColumnWidth = 23

# Restore the (now-scrolled) selection left by the author after editing.
[void]$ws.Range("J22").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 3
